$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Niam Vamyeejkoob"
$ws.Range("A3").Value = "Vamyeejkoob"

$ws.Range("A4").Select()
